# Generate Report for Handback
#
# Updates the localization-status workbook with the results of a failed
# handback transform for the file
# "1f506205-859e-43a2-96e6-c7b84ee9f2f5.103b4956155b96f72329c8335e5a1705975545ad"
# for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1. Update every cell that previously showed the "Ready for handoff"
#    status for the 1f506205... row (row 3) to reflect the failed handback
#    transform: the Overview sheet's per-locale status columns (E = zh-cn,
#    F = de-de) and each locale sheet's own Status column (C).
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# 2. zh-cn sheet: widen the "Error Detail" column (P) and record the
#    handback/handoff file-name mismatch for the 1f506205... row (row 3).
$zhcn.Range("P1").ColumnWidth = 39.17
$zhcn.Range("P3").Value = "Handback file name: h2bfn4io.vuy is different with handoff file name: 1f506205-859e-43a2-96e6-c7b84ee9f2f5.103b4956155b96f72329c8335e5a1705975545ad.zh-cn."

# 3. de-de sheet: widen the "Error Detail" column (P) and record the
#    handback/handoff file-name mismatch for the 1f506205... row (row 3).
$dede.Range("P1").ColumnWidth = 39.17
$dede.Range("P3").Value = "Handback file name: h2bfn4io.vuy is different with handoff file name: 1f506205-859e-43a2-96e6-c7b84ee9f2f5.103b4956155b96f72329c8335e5a1705975545ad.de-de."
